# Fix next possible queues
#
# The "List of possible next queues" column (D) in the Next Possible Queues
# rule table on Sheet1 needs to include every queue a case can be routed to
# (comma separated), not just a single queue, and the Release row's return
# queue mapping moves down one row so the table lines up correctly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 21 (Fulfill Next Queues): next queues should be Intake + Supervisor Approval
$ws.Range("D21").Value = "Intake,Supervisor Approval"

# Row 22 (Supervisor Approval Next Queues): next queues should be Fulfill + Executive Approval
$ws.Range("D22").Value = "Fulfill,Executive Approval"

# Row 23 (Executive Approval Next Queues): next queues should be Supervisor Approval + Release
$ws.Range("D23").Value = "Supervisor Approval,Release"

# Row 24 (Release Next Queues): next queue is Executive Approval, and the
# stray "null" default-next-queue value that belonged on this row moves out
$ws.Range("D24").Value = "Executive Approval"
$ws.Range("E24").Value = ""

# Update the active selection to reflect where the editor left off
$ws.Range("E25").Select()
